# "Fruta / hortaliza, semanal" update
#
# A new weekly price record for "Cebollín" at Feria Lagunitas de Puerto Montt
# needs to be inserted as the new first row of this block (row 356), pushing
# every subsequent record down by one row (356->357, ..., 410->411).
#
# Insert a real row (so styles/format carry over and the sheet dimension
# grows from R410 to R411) and then populate it with the new record's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 356; everything from 356 downward shifts to 357+.
$ws.Rows(356).Insert()

$row = 356
$ws.Cells.Item($row, 1).Value  = 4
$ws.Cells.Item($row, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item($row, 3).Value  = "Los Lagos"
$ws.Cells.Item($row, 4).Value  = 44984
$ws.Cells.Item($row, 5).Value  = 10
$ws.Cells.Item($row, 6).Value  = 100112037
$ws.Cells.Item($row, 7).Value  = "Cebollín"
$ws.Cells.Item($row, 8).Value  = "Sin especificar"
$ws.Cells.Item($row, 9).Value  = "Primera"
$ws.Cells.Item($row, 10).Value = 70
$ws.Cells.Item($row, 11).Value = 6500
$ws.Cells.Item($row, 12).Value = 7000
$ws.Cells.Item($row, 13).Value = 6750
$ws.Cells.Item($row, 14).Value = "$/paquete 36 unidades"
$ws.Cells.Item($row, 15).Value = "Región Metropolitana"
$ws.Cells.Item($row, 16).Value = 188
$ws.Cells.Item($row, 17).Value = 36
$ws.Cells.Item($row, 18).Value = "Hortaliza"
